# Update "想去人数" (F column) figures across the 展览, 演出 and 全部类型 sheets
# to reflect the latest scrape (output generated at 456a3b4).

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 103
$ws1.Range("F4").Value = 93
$ws1.Range("F5").Value = 2625
$ws1.Range("F6").Value = 246

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F3").Value = 3

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 103
$ws4.Range("F4").Value = 93
$ws4.Range("F5").Value = 2625
$ws4.Range("F6").Value = 246
$ws4.Range("F8").Value = 3
